$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.705.68"
$ws.Range("E2").Value = "  -7.55%  "
$ws.Range("D3").Value = "2.540.55"
$ws.Range("E3").Value = "  -4.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.88"
$ws.Range("E5").Value = "  -3.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.10"
$ws.Range("E6").Value = "  -6.93%  "
$ws.Range("E7").Value = "  -3.96%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("E9").Value = "  -5.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.81"
$ws.Range("E10").Value = "  -8.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0799"
$ws.Range("E11").Value = "  -5.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.64"
$ws.Range("E12").Value = "  -5.71%  "
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "2.926.75"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").Value = "2.576.20"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.873"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.18"
$ws.Range("E17").Value = "  -5.21%  "
$ws.Range("D18").Value = "42.737.98"
$ws.Range("E18").Value = "  -7.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("E20").Value = "  -4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.54"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.29"
$ws.Range("E22").Value = "  -4.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.13"
$ws.Range("E23").Value = "  -9.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.92"
$ws.Range("E24").Value = "  -4.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "29.26"
$ws.Range("E25").Value = "  -3.24%  "
$ws.Range("E26").Value = "  -6.63%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  -5.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.05"
$ws.Range("E29").Value = "  -4.81%  "
$ws.Range("E30").Value = "  -5.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.93"
$ws.Range("E31").Value = "  -5.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.47"
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("E33").Value = "  -7.66%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.39"
$ws.Range("E35").Value = "  -9.77%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0792"
$ws.Range("E36").Value = "  -6.24%  "
$ws.Range("E37").Value = "  -7.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.120"
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.05"
$ws.Range("E39").Value = "  +3.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.92"
$ws.Range("E40").Value = "  +6.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0310"
$ws.Range("E41").Value = "  -6.08%  "
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("E43").Value = "  -5.58%  "
$ws.Range("D44").Value = "2.082.60"
$ws.Range("E44").Value = "  -2.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.09"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.51"
$ws.Range("E47").Value = "  -10.55%  "
$ws.Range("E48").Value = "  +2.52%  "
$ws.Range("D49").Value = "2.784.05"
$ws.Range("E49").Value = "  -3.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.02"
$ws.Range("E50").Value = "  -6.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  -5.56%  "
